$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates driven by the diff. Values that look like plain
# numbers (single decimal point, e.g. "242.23") are written with a
# leading apostrophe to force text storage (matching the source
# inline-string cells), then the style is reset to Normal so no stray
# number-format/quote-prefix styling is left on the cell.

$ws.Range("D2").Value = '36.499.51'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = '2.054.51'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = "'242.23"
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'54.54"
$ws.Range("E8").Value = '  -5.62%  '
$ws.Range("D9").Value = "'58.35"
$ws.Range("E9").Value = '  -2.41%  '
$ws.Range("D10").Value = "'0.360"
$ws.Range("E10").Value = '  -4.41%  '
$ws.Range("D11").Value = "'0.0750"
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -3.05%  '
$ws.Range("D13").Value = "'0.907"
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("D14").Value = "'14.72"
$ws.Range("E14").Value = '  -4.71%  '
$ws.Range("D15").Value = '2.351.57'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = "'5.40"
$ws.Range("E16").Value = '  -4.04%  '
$ws.Range("D17").Value = '2.045.37'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '36.382.27'
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("D19").Value = "'16.76"
$ws.Range("E19").Value = '  -7.14%  '
$ws.Range("D20").Value = "'71.85"
$ws.Range("E20").Value = '  -2.51%  '
$ws.Range("D21").Value = '0.0₃0857'
$ws.Range("E21").Value = '  -3.27%  '
$ws.Range("D22").Value = "'238.42"
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").Value = "'5.24"
$ws.Range("E23").Value = '  -2.58%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = '  -4.00%  '
$ws.Range("D26").Value = "'9.34"
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = "'164.70"
$ws.Range("E28").Value = '  -2.75%  '
$ws.Range("D29").Value = "'20.05"
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("E30").Value = '  -0.99%  '
$ws.Range("D31").Value = "'1.22"
$ws.Range("E31").Value = '  +10.87%  '
$ws.Range("D32").Value = "'5.07"
$ws.Range("E32").Value = '  -4.81%  '
$ws.Range("D33").Value = "'4.45"
$ws.Range("E33").Value = '  -4.79%  '
$ws.Range("D34").Value = "'0.0594"
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = "'2.21"
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("D38").Value = "'0.0817"
$ws.Range("E38").Value = '  -6.18%  '
$ws.Range("D39").Value = "'1.24"
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("D40").Value = "'4.82"
$ws.Range("E40").Value = '  -5.04%  '
$ws.Range("D41").Value = "'0.0215"
$ws.Range("E41").Value = '  -2.95%  '
$ws.Range("D42").Value = "'0.0939"
$ws.Range("E42").Value = '  -4.25%  '
$ws.Range("E43").Value = '  -9.21%  '
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("D45").Value = "'93.74"
$ws.Range("E45").Value = '  -3.02%  '
$ws.Range("D46").Value = '1.410.21'
$ws.Range("E46").Value = '  +9.43%  '
$ws.Range("D47").Value = "'7.57"
$ws.Range("E47").Value = '  +12.82%  '
$ws.Range("D48").Value = "'15.90"
$ws.Range("E48").Value = '  -5.52%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'2.26"
$ws.Range("E50").Value = '  -2.37%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.238.23'
$ws.Range("E51").Value = '  +0.51%  '

# Reset styles on the forced-text cells back to Normal so no explicit
# quotePrefix/number-format style lingers on them.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
